$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2-5
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 3).Value = "s__Ruminococcus_F champanellensis"
    $ws.Cells.Item($r, 5).Value = 0.5
}
